$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ETH")
$ws.Range("B12").Value = 0.00747298
$ws.Range("B36").Value = 0.02570108
$ws.Range("B40").Value = 0.05826808
$ws.Range("D36").Value = 45.9
$ws.Range("D40").Value = 111.05
$ws.Range("J3").Value = 3465.126543523869

$ws = $wb.Worksheets.Item("BTC")
$ws.Range("B24").Value = 0.00168253
$ws.Range("B34").Value = 0.00218158
$ws.Range("B6").Value = 0.00035702
$ws.Range("D24").Value = 45.9
$ws.Range("D34").Value = 66.95
$ws.Range("J3").Value = 62814.06315777063

$ws = $wb.Worksheets.Item("POLIS")
$ws.Range("J3").Value = 0.5292126049109773

$ws = $wb.Worksheets.Item("ATLAS")
$ws.Range("J3").Value = 0.00709134242434479

$ws = $wb.Worksheets.Item("ACE")
$ws.Range("B6").Value = 0.00002752
$ws.Range("J3").Value = 12.33124880267018

$ws = $wb.Worksheets.Item("ADA")
$ws.Range("B6").Value = 0.79332135
$ws.Range("B7").Value = 125.37719126
$ws.Range("D7").Value = 45.9
$ws.Range("J3").Value = 0.6949611959446909

$ws = $wb.Worksheets.Item("ALGO")
$ws.Range("B6").Value = 0.58769182
$ws.Range("J3").Value = 0.2234593819908972

$ws = $wb.Worksheets.Item("AMP")
$ws.Range("J3").Value = 0.004986551711064912

$ws = $wb.Worksheets.Item("APE")
$ws.Range("B5").Value = 17.12465969
$ws.Range("B6").Value = 0.6019319
$ws.Range("D5").Value = 45.9
$ws.Range("J3").Value = 1.978488072207177

$ws = $wb.Worksheets.Item("ATOM")
$ws.Range("B7").Value = 0.02986839
$ws.Range("J3").Value = 11.82825788248545

$ws = $wb.Worksheets.Item("AVAX")
$ws.Range("B5").Value = 2.69520065
$ws.Range("B6").Value = 0.01682194
$ws.Range("D5").Value = 45.9
$ws.Range("J3").Value = 43.68926404349773

$ws = $wb.Worksheets.Item("BNB")
$ws.Range("B10").Value = 0.00281719
$ws.Range("B12").Value = 0.1586438
$ws.Range("D12").Value = 45.9
$ws.Range("J3").Value = 407.9508642673276

$ws = $wb.Worksheets.Item("DOGE")
$ws.Range("B6").Value = 0.29448781
$ws.Range("J3").Value = 0.1338975560051328

$ws = $wb.Worksheets.Item("DOT")
$ws.Range("B5").Value = 7.93841182
$ws.Range("B6").Value = 0.08116978
$ws.Range("D5").Value = 45.9
$ws.Range("J3").Value = 8.852857171128743

$ws = $wb.Worksheets.Item("EGLD")
$ws.Range("B6").Value = 0.00300679
$ws.Range("J3").Value = 62.54162406785102

$ws = $wb.Worksheets.Item("GRT")
$ws.Range("J3").Value = 0.2882271684537856

$ws = $wb.Worksheets.Item("ICP")
$ws.Range("B6").Value = 0.00237353
$ws.Range("J3").Value = 13.14921670490126

$ws = $wb.Worksheets.Item("KAVA")
$ws.Range("J3").Value = 0.8726665209206731

$ws = $wb.Worksheets.Item("LDO")
$ws.Range("B6").Value = 0.02067706
$ws.Range("J3").Value = 3.560811693986734

$ws = $wb.Worksheets.Item("LINK")
$ws.Range("B6").Value = 0.00250722
$ws.Range("J3").Value = 20.13779851056128

$ws = $wb.Worksheets.Item("LTC")
$ws.Range("B6").Value = 0.00137003
$ws.Range("J3").Value = 82.36938507573605

$ws = $wb.Worksheets.Item("LUNA")
$ws.Range("B6").Value = 0.05870622
$ws.Range("J3").Value = 0.758027344635119

$ws = $wb.Worksheets.Item("LUNC")
$ws.Range("B18").Value = 5071.50277339
$ws.Range("J3").Value = 0.0001492261537490401

$ws = $wb.Worksheets.Item("MATIC")
$ws.Range("B6").Value = 0.33079316
$ws.Range("B7").Value = 50.35525816
$ws.Range("D7").Value = 45.9
$ws.Range("J3").Value = 1.044927566587358

$ws = $wb.Worksheets.Item("MEME")
$ws.Range("B6").Value = 0.06959824000000001
$ws.Range("J3").Value = 0.03429921451474223

$ws = $wb.Worksheets.Item("MINA")
$ws.Range("B6").Value = 0.35232281
$ws.Range("J3").Value = 1.345270873257257

$ws = $wb.Worksheets.Item("NEAR")
$ws.Range("B6").Value = 24.415873
$ws.Range("B7").Value = 0.10354829
$ws.Range("D6").Value = 45.9
$ws.Range("J3").Value = 3.99715407897387

$ws = $wb.Worksheets.Item("SEI")
$ws.Range("B6").Value = 0.07657216
$ws.Range("J3").Value = 0.9028881702589886

$ws = $wb.Worksheets.Item("SHIB")
$ws.Range("B6").Value = 284.67
$ws.Range("J3").Value = 0.00001430254960894466

$ws = $wb.Worksheets.Item("SHPING")
$ws.Range("J3").Value = 0.006614258948515647

$ws = $wb.Worksheets.Item("SOL")
$ws.Range("B17").Value = 0.06511594
$ws.Range("B18").Value = 1.93147307
$ws.Range("D18").Value = 45.9
$ws.Range("J3").Value = 132.1996353159481

$ws = $wb.Worksheets.Item("TRX")
$ws.Range("B6").Value = 0.26939611
$ws.Range("J3").Value = 0.142242751465132

$ws = $wb.Worksheets.Item("UNI")
$ws.Range("B6").Value = 0.00278346
$ws.Range("J3").Value = 10.90314685515259

$ws = $wb.Worksheets.Item("XRP")
$ws.Range("B6").Value = 0.87798541
$ws.Range("J3").Value = 0.6064029195142571

$ws = $wb.Worksheets.Item("TIA")
$ws.Range("B6").Value = 0.00480053
$ws.Range("J3").Value = 17.56934253941745

$ws = $wb.Worksheets.Item("DYDX")
$ws.Range("B6").Value = 0.00100902
$ws.Range("J3").Value = 3.465775556480354
